# edit.ps1 - apply "revise context on lab04-ppt" changes
# 1. Update cached datetimeFigureOut field text (5/16/2022 -> 5/17/2022)
#    on the slide master and every slide layout's Date Placeholder.
# 2. Nudge the footer-link textbox on slide 10 down slightly.
# 3. Rewrite + reposition/resize the "seed" caption textbox on slide 13,
#    dropping the old hyperlinked "srand - C++ Reference" run in favor of
#    a plain-text explanation + bare URL on its own line.
# 4. Rewrite + reposition/resize the "bucket" caption textbox on slide 14.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text on Master + all Custom Layouts
# ---------------------------------------------------------------------
function Update-DatePlaceholders($shapes, $newText) {
  foreach ($shp in $shapes) {
    $isDate = $false
    try {
      if ($shp.PlaceholderFormat.Type -eq 16) { $isDate = $true }
    } catch {
    }
    if ($isDate) {
      $shp.TextFrame.TextRange.Text = $newText
    }
  }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes "5/17/2022"

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
  $layout = $master.CustomLayouts.Item($i)
  Update-DatePlaceholders $layout.Shapes "5/17/2022"
}

# ---------------------------------------------------------------------
# 2) Slide 10 - "文字方塊 3" footer link textbox, move down slightly
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$linkBox10 = $s10.Shapes.Item(3)
$linkBox10.Top = 509.4726867675781

# ---------------------------------------------------------------------
# 3) Slide 13 - "文字方塊 2" seed explanation textbox
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$seedBox = $s13.Shapes.Item(4)
$seedBox.Left = 0
$seedBox.Top = 489.1078186035156
$seedBox.Width = 910.5
$seedBox.Height = 50.892208099365234
$seedBox.TextFrame.TextRange.Text = "The seed is like a random table. The initialization with different seeds will generate different random tables.`rhttps://www.cplusplus.com/reference/cstdlib/srand/"

# ---------------------------------------------------------------------
# 4) Slide 14 - "文字方塊 2" bucket explanation textbox
# ---------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$bucketBox = $s14.Shapes.Item(7)
$bucketBox.Left = 0
$bucketBox.Top = 491.3578186035156
$bucketBox.Width = 942.75
$bucketBox.Height = 50.892208099365234
$bucketBox.TextFrame.TextRange.Text = "The purpose of the bucket is to make a uniform distribution between [10, 14]. Therefore, r will be discarded when r is equal to 10 or 11."
